$d = $word.ActiveDocument
$q = [char]0x2019   # right single quotation mark used by the source doc

# ---------------------------------------------------------------------
# Edit 1: the "P_VALID" row description gains an explanatory sentence
# after "...valid and usable" ("D" stays its own run; only the
# "etermines..." run and trailing new text are touched).
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "etermines if the prediction output is valid and usable",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Edit 1: could not find the target sentence"
}
$rng1.Collapse(0)
$rng1.InsertAfter(". The same as fw_valid_i delay matched with prediction result output path")
Write-Output "Edit 1 applied"

# ---------------------------------------------------------------------
# Edit 2: "the table's output" -> "the tables' output" (pluralize the
# possessive) in the predictor bullet point.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$find2 = "table" + $q + "s"
$found2 = $rng2.Find.Execute(
    $find2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Edit 2: could not find 'table's' text"
}
$rng2.Text = ""
$rng2.Collapse(0)
$rng2.InsertAfter("table")
$rng2.Collapse(0)
$rng2.InsertAfter("s")
$rng2.Collapse(0)
$rng2.InsertAfter($q)
Write-Output "Edit 2 applied"
